$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COMT")
$ws.Range("B138").Value = "You are in the Moderate Loss of Function category.  See below for more information."
